$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the table row whose first-column label is "FI" (the Finland row),
# which currently holds duplicated/incorrect placeholder values ("-") that
# need to be cleared out, same as the other fixed-up institutional factor
# rows in this table.
$targetRow = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $labelCell = $t.Cell($r, 1)
    $labelRng = $labelCell.Range
    $labelRng.End = $labelRng.End - 1
    if ($labelRng.Text -eq "FI") {
        $targetRow = $r
        break
    }
}

if ($targetRow -ne -1) {
    $colCount = $t.Columns.Count
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $t.Cell($targetRow, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $oldText = $rng.Text
        if ($oldText -ne "") {
            $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
        }
    }
}
